# Insert a new weekly data row at row 334 (pushing existing rows 334-356
# down to 335-357) and populate it with the new price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(334).Insert()

$ws.Range("A334").Value = 4
$ws.Range("B334").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C334").Value = "Los Lagos"
$ws.Range("D334").Value = 44826
$ws.Range("E334").Value = 10
$ws.Range("F334").Value = 100114014
$ws.Range("G334").Value = "Betarraga"
$ws.Range("H334").Value = "Sin especificar"
$ws.Range("I334").Value = "Primera"
$ws.Range("J334").Value = 500
$ws.Range("K334").Value = 1500
$ws.Range("L334").Value = 1500
$ws.Range("M334").Value = 1500
$ws.Range("N334").Value = "$/paquete 5 unidades"
$ws.Range("O334").Value = "Región del Maule"
$ws.Range("P334").Value = 300
$ws.Range("Q334").Value = 5
$ws.Range("R334").Value = "Hortaliza"
